# Insert a new "tags" column into the "loadbalancers" sheet's table,
# between the existing "resource_group" and "create_timeout" columns.
#
# Before: A=*name B=*subnets C=type D=resource_group E=create_timeout F=delete_timeout
# After:  A=*name B=*subnets C=type D=resource_group E=tags F=create_timeout G=delete_timeout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loadbalancers")
$lo = $ws.ListObjects.Item("Table3")

# Shift the two rightmost data columns over by one to make room for the
# new "tags" column at E, preserving their values/styles:
#   old F (delete_timeout) -> new G
#   old E (create_timeout) -> new F
$ws.Range("F1:F3").Copy($ws.Range("G1:G3"))
$ws.Range("E1:E2").Copy($ws.Range("F1:F2"))
$ws.Range("F3").Clear()

# Populate the new "tags" column header (column E); it has no data rows.
$ws.Range("E1").Value = "tags"
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

# Grow the table to include the new column G so the table definition
# (ref, autoFilter, tableColumns) reflects the new layout.
$lo.Resize($ws.Range("A1:G3"))

# Make sure the header text for the shifted columns is correct.
$ws.Range("F1").Value = "create_timeout"
$ws.Range("G1").Value = "delete_timeout"
